# Correction des fichiers exemple de C (seed + fish)
$wb = $excel.ActiveWorkbook

# --- Diet sheet: add two new rows (seed input correction) ---
$dietWs = $wb.Worksheets.Item("Diet")

$dietWs.Range("A35").Value = "Bio_diet"
$dietWs.Range("B35").Value = 1
$dietWs.Range("C35").Value = "Forage"

$dietWs.Range("A36").Value = "Bio_diet"
$dietWs.Range("B36").Value = 0.5
$dietWs.Range("C36").Value = "Wheat grain"

# --- New sheet: Energy power (fish / bioraffinery energy data) ---
$energyWs = $wb.Worksheets.Add()
$energyWs.Name = "Energy power"

# Filled column-by-column (matches the original authoring order so the
# shared-string table indices line up with the source workbook)
$energyWs.Range("A1").Value = "Facility"
$energyWs.Range("B1").Value = "Items"

$energyWs.Range("A2").Value = "Methanizer"
$energyWs.Range("A3").Value = "Methanizer"
$energyWs.Range("A4").Value = "Methanizer"
$energyWs.Range("A5").Value = "Methanizer"
$energyWs.Range("A6").Value = "Bioraffinery"
$energyWs.Range("A7").Value = "Bioraffinery"

$energyWs.Range("B2").Value = "Barley grain, Wheat grain, Oats grain, Maize corn"
$energyWs.Range("B3").Value = "waste"
$energyWs.Range("B4").Value = "bovines slurry, bovines manure, porcines slurry, porcines manure"
$energyWs.Range("B5").Value = "Wheat grain"
$energyWs.Range("B6").Value = "Forage"
$energyWs.Range("B7").Value = "Wheat grain"

$energyWs.Range("C1").Value = "Energy Power (MWh/tFW)"
$energyWs.Range("C2").Value = 0.78
$energyWs.Range("C3").Value = 0.5
$energyWs.Range("C4").Value = 0.23
$energyWs.Range("C5").Value = 0.25
$energyWs.Range("C6").Value = 0.1
$energyWs.Range("C7").Value = 0.1

$energyWs.Columns.Item(2).ColumnWidth = 57.36328125
$energyWs.Columns.Item(3).ColumnWidth = 23.81640625

# Move "Energy power" sheet to be the last tab (after Diet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$energyWs.Move($null, $lastSheet)

# Re-fetch sheet references by name since Move() invalidates old handles
$inputWs = $wb.Worksheets.Item("Input data")
$dietWs = $wb.Worksheets.Item("Diet")
$energyWs = $wb.Worksheets.Item("Energy power")

# --- View state updates ---

# "Input data" sheet: move selection from C97 to E38, leave frozen pane as-is
$inputWs.Activate()
$inputWs.Range("E38").Select()

# "Energy power" sheet: selection sits just below the entered data (B8)
$energyWs.Activate()
$energyWs.Range("B8").Select()

# "Diet" sheet becomes the active/selected tab, scrolled to row 24,
# with the cursor on the newly added data (D36)
$dietWs.Activate()
$excel.ActiveWindow.ScrollRow = 24
$dietWs.Range("D36").Select()
